$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the "last updated" footnote text (shared string reused in place).
# ---------------------------------------------------------------------------
$ws.Range("B39").Value = "Actualización: Agosto 2025."

# ---------------------------------------------------------------------------
# 2) The small Y:AA "Nacional" lookup table (rows 6-24) gains a new month
#    (Agosto 2025) at the top; every existing row shifts down by one and a
#    new row 25 appears at the bottom. Column formatting (banded style)
#    shifts down with the rows, so snapshot the original formats into a
#    scratch area first (rows offset by +1000) before overwriting anything.
# ---------------------------------------------------------------------------
for ($r = 6; $r -le 24; $r++) {
    $ws.Range("Y" + $r + ":AA" + $r).Copy()
    $scratch = $r + 1000
    $ws.Range("Y" + $scratch + ":AA" + $scratch).PasteSpecial(-4122)
}

# Row 6 takes the format that (original) row 7 had.
$ws.Range("Y1007:AA1007").Copy()
$ws.Range("Y6:AA6").PasteSpecial(-4122)

# Rows 7-25 take the format that (original) row (r-1) had.
for ($r = 25; $r -ge 7; $r--) {
    $src = ($r - 1) + 1000
    $ws.Range("Y" + $src + ":AA" + $src).Copy()
    $ws.Range("Y" + $r + ":AA" + $r).PasteSpecial(-4122)
}

# Clear the scratch area now that every format has been redistributed.
for ($r = 6; $r -le 24; $r++) {
    $scratch = $r + 1000
    $ws.Range("Y" + $scratch + ":AA" + $scratch).Clear()
}

# ---------------------------------------------------------------------------
# 3) Write the final literal values (year, month label, national rate) for
#    every row of the shifted table - row 6 is the newly-added Agosto 2025
#    figure, rows 7-25 are the previous rows' data shifted down by one.
# ---------------------------------------------------------------------------
$data = @(
    @(6, 2025, "Ago.", 2.77502412318),
    @(7, 2025, "Jul.", 2.621775430923),
    @(8, 2025, "Jun.", 2.610736839442),
    @(9, 2025, "May.", 2.282067080445),
    @(10, 2025, "Abr.", 2.454002599445),
    @(11, 2025, "Mar.", 2.85302401955),
    @(12, 2025, "Feb.", 2.613759785893),
    @(13, 2025, "Ene.", 2.7127082932),
    @(14, 2024, "Dic.", 2.745996146311),
    @(15, 2024, "Nov.", 2.877421090974),
    @(16, 2024, "Oct.", 2.95534598574),
    @(17, 2024, "Sep.", 3.129869799333),
    @(18, 2024, "Ago.", 2.652032920305),
    @(19, 2024, "Jul.", 2.931125436954),
    @(20, 2024, "Jun.", 2.822798677808),
    @(21, 2024, "May.", 2.392540680773),
    @(22, 2024, "Abr.", 2.7204319360599998),
    @(23, 2024, "Mar.", 2.999661292614),
    @(24, 2024, "Feb.", 2.454002599445),
    @(25, 2024, "Ene.", 2.85302401955)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("Y" + $r).Value = $row[1]
    $ws.Range("Z" + $r).Value = $row[2]
    $ws.Range("AA" + $r).Value = $row[3]
}

# ---------------------------------------------------------------------------
# 4) The empty Y:AA marker cells that used to sit on row 43 move down to
#    row 44 (row 43 no longer has any content at all).
# ---------------------------------------------------------------------------
$ws.Range("V44").Copy()
$ws.Range("Y44:AA44").PasteSpecial(-4122)
$ws.Range("Y43:AA43").Clear()
